$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 136, pushing the existing rows
# 136..187 down to 137..188 (dimension grows from A1:T187 to A1:T188).
$ws.Rows.Item(136).Insert()

# Populate the new row 136 with the new record's data. Columns A,B,C,E,F,
# G,H,I,J are constant across every data row in this sheet.
$ws.Cells.Item(136, 1).Value2 = 1
$ws.Cells.Item(136, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(136, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(136, 4).Value2 = 44588
$ws.Cells.Item(136, 5).Value2 = 15
$ws.Cells.Item(136, 6).Value2 = "Fruta"
$ws.Cells.Item(136, 7).Value2 = 100102
$ws.Cells.Item(136, 8).Value2 = "Cítricos"
$ws.Cells.Item(136, 9).Value2 = 100102003
$ws.Cells.Item(136, 10).Value2 = "Limón"
$ws.Cells.Item(136, 11).Value2 = "Sin especificar"
$ws.Cells.Item(136, 12).Value2 = "2a amarillo"
$ws.Cells.Item(136, 13).Value2 = 300
$ws.Cells.Item(136, 14).Value2 = 24000
$ws.Cells.Item(136, 15).Value2 = 25000
$ws.Cells.Item(136, 16).Value2 = 24500
$ws.Cells.Item(136, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(136, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item(136, 19).Value2 = 1225
$ws.Cells.Item(136, 20).Value2 = 20
